$d = $word.ActiveDocument

# 1) Remove the stray "_GoBack" bookmark wrapping the "Project" label in the
#    first table cell, while preserving the paragraph/run text.
$tbl = $d.Tables(1)
$cell = $tbl.Cell(1, 1)
$cell.Range.Delete()
$cell.Range.InsertParagraphAfter()
$cell.Range.Text = "Project"

# 2) Drop the now-incorrect ".5mA" clause from the sensor-output sentence.
$d.Content.Find.Execute(
    "exceed an output greater than 3.3v or .5mA as to not damage",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "exceed an output greater than 3.3v as to not damage", 2)

# 3) Change the PCB sensor-hat deadline from Week 7 to Week 6, and re-insert
#    a "_GoBack" bookmark right after the new "Week 6" text (splitting the
#    surrounding run the same way Word does on a manual edit).
$d.Content.Find.Execute(
    "Completing the PCB sensor hat by Week 7 as outlined",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Completing the PCB sensor hat by Week 6 as outlined", 2)

$rng = $d.Content
$rng.Find.Execute("Week 6", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $rng)
